# Experiment record.xlsx — update the estimation of error, mx and bootstrap
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header row: split "hy_stat" into three specific stat columns -------
# (write J before I so the shared-string table is built in the same order
# as the authored file: hy_stat(coverage), hy_stat(mse), hy_stat(error))
$ws.Range("J1").Value = "hy_stat(coverage)"
$ws.Range("I1").Value = "hy_stat(mse)"
$ws.Range("K1").Value = "hy_stat(error)"

# --- Mark the "Experiment" column (B) with a red fill for rows 2-7 -----
$ws.Range("B2:B7").Interior.Color = 255

# --- New experiment-3 data block (rows 8-10) ----------------------------
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 32
$ws.Range("E8").Value = 0.2
$ws.Range("F8").Value = 1.761406
$ws.Range("G8").Value = 0.3
$ws.Range("H8").Value = 0.9
$ws.Range("I8").Value = 1.6191660000000001
$ws.Range("J8").Value = 0.85
$ws.Range("K8").Value = 0.1

$ws.Range("C9").Value = 64
$ws.Range("E9").Value = 0.2
$ws.Range("F9").Value = 0.94569910000000001
$ws.Range("G9").Value = 0.3
$ws.Range("H9").Value = 0.5
$ws.Range("I9").Value = 0.67454749999999997
$ws.Range("J9").Value = 0.87
$ws.Range("K9").Value = 0.08

$ws.Range("C10").Value = 128
$ws.Range("E10").Value = 0.2
$ws.Range("F10").Value = 0.42375800000000002
$ws.Range("G10").Value = 0.2
$ws.Range("H10").Value = 0.9
$ws.Range("I10").Value = 0.38481460000000001
$ws.Range("J10").Value = 0.91
$ws.Range("K10").Value = 0.04

# --- Stray formatted (wrap-text) cell a few rows down -------------------
$ws.Range("F13").WrapText = $true

# --- Column widths for the new I/J/K columns ----------------------------
$ws.Columns.Item(9).ColumnWidth = 10.8
$ws.Range("J1:K1").EntireColumn.ColumnWidth = 14.36

# --- View: zoom to 115%, select G8, drop the old top-left-cell freeze ---
$excel.ActiveWindow.Zoom = 115
$ws.Range("G8").Select()
